$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (per-question marks for right/wrong answers)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Pull the counts (row 10) and the (now updated) marking scheme (row 11)
$right = $ws.Range("B10").Value()
$wrong = $ws.Range("C10").Value()
$maxMarks = $ws.Range("E10").Value()

$rightMark = $ws.Range("B11").Value()
$wrongMark = $ws.Range("C11").Value()

# Raw (full floating point precision) products - used for the fraction text below
$totalRightRaw = $right * $rightMark
$totalWrongRaw = $wrong * $wrongMark

# Cell values are stored using the "clean" shortest round-trip representation,
# same as Excel/.NET's default double formatting.
$totalRight = [double]($totalRightRaw.ToString())
$totalWrong = [double]($totalWrongRaw.ToString())

$ws.Range("B12").Value = $totalRight
$ws.Range("C12").Value = $totalWrong

# The displayed fraction keeps the un-rounded (full precision) numerator.
$scoreRaw = $totalRightRaw + $totalWrongRaw + 0
$maxScore = $maxMarks * $rightMark

$scoreText = $scoreRaw.ToString("G17")
$maxScoreText = $maxScore.ToString("G17")

$ws.Range("E12").Value = "$scoreText/$maxScoreText"
